$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the new "Handheld" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "outsource"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Handheld"

# --- Populate the new "Handheld" sheet with the break-down tasks ---
$values = @(
    "GET api/v2/{siteId}/che/{cheId}/workinstructions",
    "GET api/v2/{siteId}/req-wis",
    "GET api/v2/{siteId}/res-wis/{fromDate}/{toDate}",
    "GET api/v2/{siteId}/done-wis/{fromDate}/{toDate}",
    "GET api/v2/{siteId}/del-req-wis",
    "GET api/v2/{siteId}/del-res-wis/{fromData}/{toDate}",
    "PUT api/v2/{siteId}/che/{cheId}/wis/{wiid}"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $values[$i]
}

$ws2.Columns("A").ColumnWidth = 85.140625

$ws2.Range("H16").Select()

$wb.Save()
